$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 28, shifting existing rows 28.. down by one.
$ws.Rows.Item(28).Insert()

# Populate the new row 28 with the new data record.
$ws.Cells.Item(28, 1).Value = 11
$ws.Cells.Item(28, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(28, 3).Value = "Bíobío"
$ws.Cells.Item(28, 4).Value = 44707
$ws.Cells.Item(28, 5).Value = 8
$ws.Cells.Item(28, 6).Value = 100112043
$ws.Cells.Item(28, 7).Value = "Pepino ensalada"
$ws.Cells.Item(28, 8).Value = "Sin especificar"
$ws.Cells.Item(28, 9).Value = "Primera"
$ws.Cells.Item(28, 10).Value = 100
$ws.Cells.Item(28, 11).Value = 18000
$ws.Cells.Item(28, 12).Value = 20000
$ws.Cells.Item(28, 13).Value = 19000
$ws.Cells.Item(28, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(28, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(28, 16).Value = 317
$ws.Cells.Item(28, 17).Value = 60
$ws.Cells.Item(28, 18).Value = "Hortaliza"

# Apply the same date style (numFmt) as the other rows in column D.
$ws.Cells.Item(28, 4).NumberFormat = $ws.Cells.Item(29, 4).NumberFormat
